$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Delete the entire "VidStreamDownloadBtnPostDelay" row (row 11), shifting
# everything below it up by one row.
$ws.Rows.Item(11).Delete()

# Update the "DelayShort" value (now on row 13 after the deletion) from
# 2000 to 3000.
$ws.Range("B13").Value = 3000

# Move the active selection, mirroring the author's last click.
$ws.Range("A14").Select()
